# Applies the "Updated symbol list" refresh (Tue Feb 7 20:08:36 UTC 2023) to
# the cryptos sheet: every data row's "Hora" (G) column moves from "19" to
# "20", and most rows' Price (D) / Volume(1h) (E) columns are refreshed with
# the newly scraped figures.
#
# All of these cells are plain text cells (t="inlineStr") holding
# number-shaped strings (e.g. "0.1100", "0.0002000") and percentages
# (e.g. "-0.52%"). Assigning a numeric-looking string straight to
# Range.Value lets Excel auto-coerce it to a real number, which would
# silently strip significant trailing zeros (e.g. "0.1100" -> 0.11) and
# flip the cell's stored type away from text. To keep the cells textual and
# byte-for-byte faithful to the source strings, we briefly mark the range as
# Text ("@") before writing, then restore the cell's style afterwards so no
# stray formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText($sheet, $addr, $text) {
    $range = $sheet.Range($addr)
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

Set-CellText $ws 'D2' '329.99'
Set-CellText $ws 'G2' '20'
Set-CellText $ws 'D3' '44.03'
Set-CellText $ws 'E3' '-0.52%'
Set-CellText $ws 'G3' '20'
Set-CellText $ws 'D4' '5.527'
Set-CellText $ws 'E4' '-0.77%'
Set-CellText $ws 'G4' '20'
Set-CellText $ws 'D5' '0.08102'
Set-CellText $ws 'E5' '0.09%'
Set-CellText $ws 'G5' '20'
Set-CellText $ws 'D6' '2.066'
Set-CellText $ws 'E6' '4.31%'
Set-CellText $ws 'G6' '20'
Set-CellText $ws 'D7' '0.9728'
Set-CellText $ws 'E7' '2.03%'
Set-CellText $ws 'G7' '20'
Set-CellText $ws 'D8' '0.1100'
Set-CellText $ws 'E8' '-6.45%'
Set-CellText $ws 'G8' '20'
Set-CellText $ws 'E9' '2.00%'
Set-CellText $ws 'G9' '20'
Set-CellText $ws 'D10' '10.19'
Set-CellText $ws 'E10' '-0.60%'
Set-CellText $ws 'G10' '20'
Set-CellText $ws 'D11' '0.09987'
Set-CellText $ws 'E11' '0.82%'
Set-CellText $ws 'G11' '20'
Set-CellText $ws 'D12' '0.04729'
Set-CellText $ws 'E12' '0.29%'
Set-CellText $ws 'G12' '20'
Set-CellText $ws 'D13' '0.1056'
Set-CellText $ws 'E13' '-1.29%'
Set-CellText $ws 'G13' '20'
Set-CellText $ws 'D14' '0.001261'
Set-CellText $ws 'E14' '-1.86%'
Set-CellText $ws 'G14' '20'
Set-CellText $ws 'D15' '0.04095'
Set-CellText $ws 'E15' '-3.02%'
Set-CellText $ws 'G15' '20'
Set-CellText $ws 'D16' '0.006016'
Set-CellText $ws 'E16' '1.63%'
Set-CellText $ws 'G16' '20'
Set-CellText $ws 'D17' '3.341'
Set-CellText $ws 'E17' '-0.91%'
Set-CellText $ws 'G17' '20'
Set-CellText $ws 'E18' '2.38%'
Set-CellText $ws 'G18' '20'
Set-CellText $ws 'D19' '2.647'
Set-CellText $ws 'E19' '2.82%'
Set-CellText $ws 'G19' '20'
Set-CellText $ws 'D20' '0.3313'
Set-CellText $ws 'E20' '-4.54%'
Set-CellText $ws 'G20' '20'
Set-CellText $ws 'D21' '0.1390'
Set-CellText $ws 'G21' '20'
Set-CellText $ws 'D22' '0.2569'
Set-CellText $ws 'E22' '2.38%'
Set-CellText $ws 'G22' '20'
Set-CellText $ws 'D23' '0.001309'
Set-CellText $ws 'E23' '4.79%'
Set-CellText $ws 'G23' '20'
Set-CellText $ws 'D24' '0.004393'
Set-CellText $ws 'E24' '1.68%'
Set-CellText $ws 'G24' '20'
Set-CellText $ws 'D25' '0.0001280'
Set-CellText $ws 'G25' '20'
Set-CellText $ws 'D26' '0.0003735'
Set-CellText $ws 'E26' '-6.12%'
Set-CellText $ws 'G26' '20'
Set-CellText $ws 'G27' '20'
Set-CellText $ws 'G28' '20'
Set-CellText $ws 'G29' '20'
Set-CellText $ws 'G30' '20'
Set-CellText $ws 'G31' '20'
Set-CellText $ws 'G32' '20'
Set-CellText $ws 'G33' '20'
Set-CellText $ws 'G34' '20'
Set-CellText $ws 'G35' '20'
Set-CellText $ws 'G36' '20'
Set-CellText $ws 'G37' '20'
Set-CellText $ws 'D38' '0.02681'
Set-CellText $ws 'E38' '0.82%'
Set-CellText $ws 'G38' '20'
Set-CellText $ws 'D39' '0.05628'
Set-CellText $ws 'E39' '1.02%'
Set-CellText $ws 'G39' '20'
Set-CellText $ws 'D40' '0.007616'
Set-CellText $ws 'E40' '0.54%'
Set-CellText $ws 'G40' '20'
Set-CellText $ws 'D41' '0.1414'
Set-CellText $ws 'E41' '0.24%'
Set-CellText $ws 'G41' '20'
Set-CellText $ws 'D42' '0.007520'
Set-CellText $ws 'E42' '-6.91%'
Set-CellText $ws 'G42' '20'
Set-CellText $ws 'D43' '0.001958'
Set-CellText $ws 'G43' '20'
Set-CellText $ws 'D44' '0.008329'
Set-CellText $ws 'E44' '-6.40%'
Set-CellText $ws 'G44' '20'
Set-CellText $ws 'D45' '0.00007022'
Set-CellText $ws 'E45' '-2.72%'
Set-CellText $ws 'G45' '20'
Set-CellText $ws 'D46' '0.00000000750'
Set-CellText $ws 'E46' '-0.18%'
Set-CellText $ws 'G46' '20'
Set-CellText $ws 'D47' '0.0005792'
Set-CellText $ws 'E47' '-0.33%'
Set-CellText $ws 'G47' '20'
Set-CellText $ws 'D48' '0.002517'
Set-CellText $ws 'E48' '10.71%'
Set-CellText $ws 'G48' '20'
Set-CellText $ws 'D49' '0.003536'
Set-CellText $ws 'E49' '-25.56%'
Set-CellText $ws 'G49' '20'
Set-CellText $ws 'D50' '0.00002099'
Set-CellText $ws 'E50' '-0.18%'
Set-CellText $ws 'G50' '20'
Set-CellText $ws 'D51' '0.0002000'
Set-CellText $ws 'E51' '-0.18%'
Set-CellText $ws 'G51' '20'
